# Apply the "calculation of new indicators" edit across the workbook.
#
# Summary of changes:
#  - SCHEME_MEASURES:   MQMS01..MQMS05  -> MQME001..MQME005   (indicator codes only)
#  - METADATA_ISSUES:   MQME10 -> MQME012, MQME12 -> MQME014,
#                       MQME01 -> MQME008, MQME16 -> MQME011  (indicator codes only)
#  - METADATA_MEASURES: re-numbered indicators, descriptions/values shifted up
#                       one row, and the last (now-duplicate) row is removed.
#  - METADATA_METRICS:  re-numbered indicators, descriptions/values shifted
#                       down, and four new indicator rows are appended.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) SCHEME_MEASURES - rename indicator codes MQMS0n -> MQME00n
# ---------------------------------------------------------------------------
$wsScheme = $wb.Worksheets.Item("SCHEME_MEASURES")
$wsScheme.Range("A2").Value = "MQME001"
$wsScheme.Range("A3").Value = "MQME002"
$wsScheme.Range("A4").Value = "MQME003"
$wsScheme.Range("A5").Value = "MQME004"
$wsScheme.Range("A6").Value = "MQME005"

# ---------------------------------------------------------------------------
# 2) METADATA_ISSUES - rename indicator codes in column A (rows 2-225)
#    Use whole-cell Replace to avoid partial / overlapping substring matches
#    (e.g. "MQME01" is a substring of the new code "MQME012").
# ---------------------------------------------------------------------------
$wsIssues = $wb.Worksheets.Item("METADATA_ISSUES")
$issuesRange = $wsIssues.Range("A2:A225")
$issuesRange.Replace("MQME10", "MQME012", 1)
$issuesRange.Replace("MQME12", "MQME014", 1)
$issuesRange.Replace("MQME01", "MQME008", 1)
$issuesRange.Replace("MQME16", "MQME011", 1)

# ---------------------------------------------------------------------------
# 3) METADATA_MEASURES - renumber & shift rows, drop the last row
# ---------------------------------------------------------------------------
$wsMeasures = $wb.Worksheets.Item("METADATA_MEASURES")
$wsMeasures.Range("A2").Value = "MQME006"
$wsMeasures.Range("B2").Value = "Total number of length-required columns"
$wsMeasures.Range("C2").Value = 179
$wsMeasures.Range("A3").Value = "MQME007"
$wsMeasures.Range("B3").Value = "Total number of NUMBER columns"
$wsMeasures.Range("C3").Value = 363
$wsMeasures.Rows.Item(4).Delete()

# ---------------------------------------------------------------------------
# 4) METADATA_METRICS - renumber & shift rows, append four new rows
# ---------------------------------------------------------------------------
$wsMetrics = $wb.Worksheets.Item("METADATA_METRICS")

function Set-TextValue($ws, $addr, $text) {
    # Percentage-looking strings (e.g. "90.91%") get auto-converted to a
    # numeric percentage by Excel unless the cell is explicitly formatted as
    # text first. Reset the style back to Normal afterwards so no stray
    # cell formatting is introduced.
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

Set-TextValue $wsMetrics "A2" "MQID001"
Set-TextValue $wsMetrics "B2" "Table names in singular"
Set-TextValue $wsMetrics "C2" "95.12%"

Set-TextValue $wsMetrics "A3" "MQID002"
Set-TextValue $wsMetrics "B3" "Table with recommended name length"
Set-TextValue $wsMetrics "C3" "100.00%"

Set-TextValue $wsMetrics "A4" "MQID003"
Set-TextValue $wsMetrics "B4" "Columns with correct prefixes"
Set-TextValue $wsMetrics "C4" "99.32%"

Set-TextValue $wsMetrics "A5" "MQID004"
Set-TextValue $wsMetrics "B5" "Columns with recommended name size"
Set-TextValue $wsMetrics "C5" "100.00%"

Set-TextValue $wsMetrics "A6" "MQID005"
Set-TextValue $wsMetrics "B6" "Columns with comments"
Set-TextValue $wsMetrics "C6" "63.16%"

Set-TextValue $wsMetrics "A7" "MQID006"
Set-TextValue $wsMetrics "B7" "Table with standard PK prefixes"
Set-TextValue $wsMetrics "C7" "100.00%"

Set-TextValue $wsMetrics "A8" "MQID007"
Set-TextValue $wsMetrics "B8" "Table with standard FK prefixes"
Set-TextValue $wsMetrics "C8" "100.00%"

Set-TextValue $wsMetrics "A9" "MQID008"
Set-TextValue $wsMetrics "B9" "Table with standard UK prefixes"
Set-TextValue $wsMetrics "C9" "90.91%"

Set-TextValue $wsMetrics "A10" "MQID009"
Set-TextValue $wsMetrics "B10" "NUMBER columns with valid scale"
Set-TextValue $wsMetrics "C10" "100.00%"

Set-TextValue $wsMetrics "A11" "MQID010"
Set-TextValue $wsMetrics "B11" "Columns with valid num_distinct"
Set-TextValue $wsMetrics "C11" "100.00%"

Set-TextValue $wsMetrics "A12" "MQID011"
Set-TextValue $wsMetrics "B12" "Columns with valid num_nulls"
Set-TextValue $wsMetrics "C12" "100.00%"
